# Update res_bus q_mvar results after quadrupling the wind power input.
# Recalculated reactive power (Q, MVAr) values for buses 0-23 (rows 2-25)
# in columns B, F, I, J, K.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 46.16748940406251
$ws.Range("F2").Value = 110.0290523767471
$ws.Range("I2").Value = 51.36856651306152
$ws.Range("J2").Value = 23.17655599117279
$ws.Range("K2").Value = 12.89958465099335
$ws.Range("B3").Value = 43.60261143407661
$ws.Range("F3").Value = 110.4455714225769
$ws.Range("I3").Value = 52.08560562133789
$ws.Range("J3").Value = 24.11791467666626
$ws.Range("K3").Value = 14.67134952545166
$ws.Range("B4").Value = 43.95248131707376
$ws.Range("F4").Value = 110.830887556076
$ws.Range("I4").Value = 52.10620939731598
$ws.Range("J4").Value = 24.90840172767639
$ws.Range("K4").Value = 14.36756181716919
$ws.Range("B5").Value = 44.54332154912845
$ws.Range("F5").Value = 111.1066139936447
$ws.Range("I5").Value = 51.44108867645264
$ws.Range("J5").Value = 25.35288536548615
$ws.Range("K5").Value = 14.56564688682556
$ws.Range("B6").Value = 44.82914441899629
$ws.Range("F6").Value = 111.2533802986145
$ws.Range("I6").Value = 50.89379394054413
$ws.Range("J6").Value = 25.55463826656342
$ws.Range("K6").Value = 14.91156375408173
$ws.Range("B7").Value = 43.03911501550829
$ws.Range("F7").Value = 111.2521994113922
$ws.Range("I7").Value = 50.88704538345337
$ws.Range("J7").Value = 25.52816152572632
$ws.Range("K7").Value = 16.62175238132477
$ws.Range("B8").Value = 42.67879587397874
$ws.Range("F8").Value = 111.2446665763855
$ws.Range("I8").Value = 50.74931299686432
$ws.Range("J8").Value = 25.48218643665314
$ws.Range("K8").Value = 17.11912178993225
$ws.Range("B9").Value = 43.79046449522684
$ws.Range("F9").Value = 110.980545759201
$ws.Range("I9").Value = 50.29239797592163
$ws.Range("J9").Value = 24.87215793132782
$ws.Range("K9").Value = 16.57189047336578
$ws.Range("B10").Value = 47.37179880942494
$ws.Range("F10").Value = 110.1932618618011
$ws.Range("I10").Value = 50.18911325931549
$ws.Range("J10").Value = 23.30446016788483
$ws.Range("K10").Value = 13.09374070167542
$ws.Range("B11").Value = 51.48168642928795
$ws.Range("F11").Value = 109.4196938276291
$ws.Range("I11").Value = 49.9255610704422
$ws.Range("J11").Value = 21.76256930828094
$ws.Range("K11").Value = 9.044375658035278
$ws.Range("B12").Value = 53.10579382179822
$ws.Range("F12").Value = 109.078875541687
$ws.Range("I12").Value = 49.41969001293182
$ws.Range("J12").Value = 21.17377579212189
$ws.Range("K12").Value = 7.637016534805298
$ws.Range("B13").Value = 53.50864612417354
$ws.Range("F13").Value = 108.8093898296356
$ws.Range("I13").Value = 48.97060418128967
$ws.Range("J13").Value = 21.61095023155212
$ws.Range("K13").Value = 6.530222773551941
$ws.Range("B14").Value = 53.12059195717484
$ws.Range("F14").Value = 108.782541513443
$ws.Range("I14").Value = 48.64319276809692
$ws.Range("J14").Value = 22.42372572422028
$ws.Range("K14").Value = 6.346652984619141
$ws.Range("B15").Value = 52.92441167664765
$ws.Range("F15").Value = 108.8951833248138
$ws.Range("I15").Value = 48.28096687793732
$ws.Range("J15").Value = 22.94952356815338
$ws.Range("K15").Value = 6.649441003799438
$ws.Range("B16").Value = 52.32623580147265
$ws.Range("F16").Value = 108.8515141010284
$ws.Range("I16").Value = 48.70133566856384
$ws.Range("J16").Value = 23.0654444694519
$ws.Range("K16").Value = 6.629059791564941
$ws.Range("B17").Value = 51.63867588068024
$ws.Range("F17").Value = 108.7909594774246
$ws.Range("I17").Value = 49.24281930923462
$ws.Range("J17").Value = 23.13435959815979
$ws.Range("K17").Value = 6.588289618492126
$ws.Range("B18").Value = 50.79171067350217
$ws.Range("F18").Value = 109.0281347036362
$ws.Range("I18").Value = 49.39900290966034
$ws.Range("J18").Value = 23.52727448940277
$ws.Range("K18").Value = 7.503948926925659
$ws.Range("B19").Value = 50.22294633708952
$ws.Range("F19").Value = 109.1468261480331
$ws.Range("I19").Value = 49.76924443244934
$ws.Range("J19").Value = 23.49616575241089
$ws.Range("K19").Value = 8.053111791610718
$ws.Range("B20").Value = 51.45509313308685
$ws.Range("F20").Value = 109.1592726707458
$ws.Range("I20").Value = 49.71204161643982
$ws.Range("J20").Value = 22.34852719306946
$ws.Range("K20").Value = 8.050723791122437
$ws.Range("B21").Value = 53.79503551165726
$ws.Range("F21").Value = 108.7651250362396
$ws.Range("I21").Value = 49.64154195785522
$ws.Range("J21").Value = 20.54062652587891
$ws.Range("K21").Value = 6.564741611480713
$ws.Range("B22").Value = 54.73844762309454
$ws.Range("F22").Value = 108.3314599990845
$ws.Range("I22").Value = 49.92917191982269
$ws.Range("J22").Value = 19.66533887386322
$ws.Range("K22").Value = 5.07741117477417
$ws.Range("B23").Value = 55.1727860473693
$ws.Range("F23").Value = 107.973271727562
$ws.Range("I23").Value = 50.34068858623505
$ws.Range("J23").Value = 19.04461634159088
$ws.Range("K23").Value = 3.910771131515503
$ws.Range("B24").Value = 55.08814650893692
$ws.Range("F24").Value = 107.815279006958
$ws.Range("I24").Value = 50.71421527862549
$ws.Range("J24").Value = 18.80286824703217
$ws.Range("K24").Value = 3.46308970451355
$ws.Range("B25").Value = 52.41374612631262
$ws.Range("F25").Value = 108.335921049118
$ws.Range("I25").Value = 51.55812072753906
$ws.Range("J25").Value = 19.93945348262787
$ws.Range("K25").Value = 5.655738711357117
